$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Update UI state on the current last sheet before the insertion point
# ("AlertsNotificationManagement") - it stops being the active/selected tab
# and its lingering cell selection moves to C16.
# ---------------------------------------------------------------------------
$anchor = $wb.Worksheets.Item("AlertsNotificationManagement")
$anchor.Activate() | Out-Null
$anchor.Range("C16").Select() | Out-Null

# ---------------------------------------------------------------------------
# Insert the new worksheet "CoreAlertsMaskAcAndMaskingAutho" right after
# "AlertsNotificationManagement" (becomes the 6th tab).
# ---------------------------------------------------------------------------
$new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $anchor)
$new.Name = "CoreAlertsMaskAcAndMaskingAutho"

# Reuse the formatting (fonts/fills/borders) of a structurally similar sheet
# ("AlertsTemplateManagement") for the header/data rows.
$template = $wb.Worksheets.Item("AlertsTemplateManagement")
$template.Range("A1:H2").Copy() | Out-Null
$new.Range("A1:H2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Header row
$new.Range("A1").Value = "TestScenario"
$new.Range("B1").Value = "Run"
$new.Range("C1").Value = "Status"
$new.Range("D1").Value = "Reason"
$new.Range("E1").Value = "NotificationType"
$new.Range("F1").Value = "Template"
$new.Range("G1").Value = "Status2"
$new.Range("H1").Value = "InvalidAccountNumber"

# Data row
$new.Range("A2").Value = "CoreAlertsMaskAcAndMaskingAutho"
$new.Range("B2").Value = "Yes"
$new.Range("C2").Value = "Mask"
$new.Range("D2").Value = "Test"
$new.Range("E2").Value = "EMAIL"
$new.Range("F2").Value = "Email 2"
$new.Range("G2").Value = "Unmask"
$new.Range("H2").Value = 1234

# Column A width to match the source workbook's formatting.
$new.Columns.Item(1).ColumnWidth = 35.8

# This new sheet becomes the active / selected tab, with a lingering
# selection on C10.
$new.Activate() | Out-Null
$new.Range("C10").Select() | Out-Null

Write-Host "CoreAlertsMaskAcAndMaskingAutho sheet added"
